$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# Update min_time/max_time values for rows 13 and 14 (item time scaling)
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 5

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 5

# Move the active selection to D14
$ws.Activate()
$ws.Range("D14").Select()
